$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: a couple of recalculated figures changed slightly, and the
# user's selection moved from A5 to C5.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A3").Value = 45.16
$wsSummary.Range("B3").Value = 45.16
$wsSummary.Range("A4").Value = 0
$wsSummary.Range("C4").Value = 0

# ---------------------------------------------------------------------------
# Repayment schedule sheet: the loan got fully settled in a single
# transaction, so the two-row tail of the schedule collapses into one row
# (old row 4 is removed, and row 3 now carries the closing totals).
# ---------------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Row 3 picks up the "amount" / "principal" number formats that row 4 used
# to carry (copy the formatting over before row 4 disappears), and two new
# cells (P2/P3) are created matching the neighbouring column's format.
$wsSchedule.Range("F4").Copy()
$wsSchedule.Range("F3").PasteSpecial(-4122)
$wsSchedule.Range("H3").Copy()
$wsSchedule.Range("G3").PasteSpecial(-4122)
$wsSchedule.Range("K4").Copy()
$wsSchedule.Range("K3").PasteSpecial(-4122)
$wsSchedule.Range("L4").Copy()
$wsSchedule.Range("L3").PasteSpecial(-4122)
$wsSchedule.Range("N2").Copy()
$wsSchedule.Range("P2").PasteSpecial(-4122)
$wsSchedule.Range("N3").Copy()
$wsSchedule.Range("P3").PasteSpecial(-4122)

$wsSchedule.Range("B3").Value = 14
$wsSchedule.Range("C3").Value = 42019
$wsSchedule.Range("F3").Value = 10000
$wsSchedule.Range("G3").Value = 0
$wsSchedule.Range("H3").Value = 45.16
$wsSchedule.Range("I3").Value = 0
$wsSchedule.Range("J3").Value = 0
$wsSchedule.Range("K3").Value = 10045.16
$wsSchedule.Range("L3").Value = 10045.16
$wsSchedule.Range("M3").Value = 0
$wsSchedule.Range("N3").Value = 0
$wsSchedule.Range("P3").Value = 100
$wsSchedule.Range("Q3").Value = 988.49

# Old row 4 (the second instalment) is now redundant - remove it entirely,
# shifting the dimension from A1:Q4 down to A1:Q3.
$wsSchedule.Rows.Item(4).Delete()

$wsSchedule.Columns.Item(11).ColumnWidth = 8.26
$wsSchedule.Columns.Item(12).ColumnWidth = 8.26

# ---------------------------------------------------------------------------
# Transactions sheet: same closing-balance figures, recalculated.
# ---------------------------------------------------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("E2").Value = 10045.16
$wsTransactions.Range("G2").Value = 45.16
$wsTransactions.Columns.Item(5).ColumnWidth = 8.26

# ---------------------------------------------------------------------------
# Prepay Loan sheet: same recalculated figure.
# ---------------------------------------------------------------------------
$wsPrepay = $wb.Worksheets.Item("Prepay Loan")
$wsPrepay.Range("B4").Value = 10045.16

# ---------------------------------------------------------------------------
# Selections / active sheet — replays the user's final navigation: they
# moved the cursor on NewLoanInput, Repayment schedule and Transactions,
# then ended up on the Prepay Loan tab (which becomes the active sheet).
# ---------------------------------------------------------------------------
$wsNewLoanInput = $wb.Worksheets.Item("NewLoanInput")
$wsNewLoanInput.Range("B9").Select()

$wsSummary.Range("C5").Select()

$wsSchedule.Range("H2").Select()

$wsTransactions.Range("E10").Select()

$wsPrepay.Range("B4").Select()
$wsPrepay.Activate()
